# remove table location llm output
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) text updates
$ws.Range("G1").Value = "table_header_position"
$ws.Range("Q1").Value = "Swaging Dimensions Swage Dia. In Mm"
$ws.Range("R1").Value = "Tail Id Top Tolerance Mm"
$ws.Range("S1").Value = "Tail Id Bottom Tolerance Mm"
$ws.Range("T1").Value = "Swaged No Go Mm"
$ws.Range("U1").Value = "Swaged Go Mm"

# Rows 5-12: "1 Wire, 100 Metre Coils" descriptor -> new text
for ($r = 5; $r -le 12; $r++) {
    $ws.Cells.Item($r, 7).Value = '- "1 Wire, 100 Metre Coils - Can be found on the bottom left position of the page"'
}

# Rows 13-23: Hydraulic Hose technical data descriptor -> new text
for ($r = 13; $r -le 23; $r++) {
    $ws.Cells.Item($r, 7).Value = '- "Technical Data - Can be found on the very bottom position of the page"'
}
